$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, pushing existing rows 48-97 down to 49-98.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new data record.
$ws.Cells.Item(48, 1).Value = 9
$ws.Cells.Item(48, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(48, 3).Value = "Metropolitana"
$ws.Cells.Item(48, 4).Value = 44651
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 13
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100101
$ws.Cells.Item(48, 8).Value = "Berries"
$ws.Cells.Item(48, 9).Value = 100101004
$ws.Cells.Item(48, 10).Value = "Frambuesa"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 380
$ws.Cells.Item(48, 14).Value = 8000
$ws.Cells.Item(48, 15).Value = 8000
$ws.Cells.Item(48, 16).Value = 8000
$ws.Cells.Item(48, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(48, 18).Value = "Provincia de Linares"
$ws.Cells.Item(48, 19).Value = 4000
$ws.Cells.Item(48, 20).Value = 2
